$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 3, shifting current rows 3-4 down to 5-6.
$ws.Rows.Item(3).Resize(2).Insert()

# --- New row 3 (ITALY - SERIE A, Lecce vs Verona) ---
$ws.Range("A3").Value = "4CLi9Djd"
$ws.Range("B3").Value = "29/10/2024"
$ws.Range("C3").Value = "14:30"
$ws.Range("D3").Value = "ITALY - SERIE A"
$ws.Range("E3").Value = "Lecce"
$ws.Range("F3").Value = "Verona"
$ws.Range("G3").Value = 2.3
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 3.4
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 4
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.67
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.5
$ws.Range("U3").Value = 1.95
$ws.Range("V3").Value = 1.8
$ws.Range("W3").Value = 6.5
$ws.Range("X3").Value = 10
$ws.Range("Y3").Value = 9.5
$ws.Range("Z3").Value = 21
$ws.Range("AA3").Value = 21
$ws.Range("AB3").Value = 34
$ws.Range("AC3").Value = 8
$ws.Range("AD3").Value = 6
$ws.Range("AE3").Value = 17
$ws.Range("AF3").Value = 51
$ws.Range("AG3").Value = 401
$ws.Range("AH3").Value = 8.5
$ws.Range("AI3").Value = 15
$ws.Range("AJ3").Value = 12
$ws.Range("AK3").Value = 34
$ws.Range("AL3").Value = 29
$ws.Range("AM3").Value = 41
$ws.Range("AN3").Value = 4.33
$ws.Range("AO3").Value = 13
$ws.Range("AP3").Value = 26
$ws.Range("AQ3").Value = 41
$ws.Range("AR3").Value = 67
$ws.Range("AS3").Value = 201
$ws.Range("AT3").Value = 2.5
$ws.Range("AU3").Value = 8.5
$ws.Range("AV3").Value = 67
$ws.Range("AW3").Value = 5
$ws.Range("AX3").Value = 19
$ws.Range("AY3").Value = 29
$ws.Range("AZ3").Value = 67
$ws.Range("BA3").Value = 101
$ws.Range("BB3").Value = 251
$ws.Range("BC3").Value = 126
$ws.Range("BD3").Value = 126

# --- New row 4 (AUSTRIA - 2. LIGA, Floridsdorfer AC vs Liefering) ---
$ws.Range("A4").Value = "IaB1Slhd"
$ws.Range("B4").Value = "29/10/2024"
$ws.Range("C4").Value = "14:30"
$ws.Range("D4").Value = "AUSTRIA - 2. LIGA"
$ws.Range("E4").Value = "Floridsdorfer AC"
$ws.Range("F4").Value = "Liefering"
$ws.Range("G4").Value = 2.57
$ws.Range("H4").Value = 3.25
$ws.Range("I4").Value = 2.6
$ws.Range("J4").Value = 3.2
$ws.Range("K4").Value = 2.07
$ws.Range("L4").Value = 3.2
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 7.4
$ws.Range("O4").Value = 1.28
$ws.Range("P4").Value = 3.35
$ws.Range("Q4").Value = 1.85
$ws.Range("R4").Value = 1.88
$ws.Range("S4").Value = 1.42
$ws.Range("T4").Value = 2.67
$ws.Range("U4").Value = 1.65
$ws.Range("V4").Value = 2.1
$ws.Range("W4").Value = 9.25
$ws.Range("X4").Value = 13.5
$ws.Range("Y4").Value = 9.5
$ws.Range("Z4").Value = 29
$ws.Range("AA4").Value = 20
$ws.Range("AB4").Value = 27
$ws.Range("AC4").Value = 7.4
$ws.Range("AD4").Value = 6.3
$ws.Range("AE4").Value = 12.5
$ws.Range("AF4").Value = 55
$ws.Range("AG4").Value = 350
$ws.Range("AH4").Value = 9
$ws.Range("AI4").Value = 13.5
$ws.Range("AJ4").Value = 9.5
$ws.Range("AK4").Value = 29
$ws.Range("AL4").Value = 21
$ws.Range("AM4").Value = 28
$ws.Range("AN4").Value = 4.55
$ws.Range("AO4").Value = 14
$ws.Range("AP4").Value = 21
$ws.Range("AQ4").Value = 60
$ws.Range("AR4").Value = 90
$ws.Range("AS4").Value = 250
$ws.Range("AT4").Value = 2.67
$ws.Range("AU4").Value = 6.9
$ws.Range("AV4").Value = 60
$ws.Range("AW4").Value = 4.55
$ws.Range("AX4").Value = 14
$ws.Range("AY4").Value = 22
$ws.Range("AZ4").Value = 65
$ws.Range("BA4").Value = 100
$ws.Range("BB4").Value = 300
$ws.Range("BC4").Value = 81
$ws.Range("BD4").Value = 81
